$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 982.625
$ws.Range("I92").Value = 238.88889
$ws.Range("K92").Value = 238.88889
$ws.Range("M92").Value = 1009.11111

$ws.Range("H94").Value = 375
$ws.Range("I94").Value = 375
$ws.Range("K94").Value = 375
$ws.Range("M94").Value = 76

$ws.Range("H100").Value = 656.8570999999999
$ws.Range("I100").Value = 599.6667
$ws.Range("K100").Value = 599.6667
$ws.Range("M100").Value = -58.66669999999999

$ws.Range("H103").Value = 900.6
$ws.Range("I103").Value = 1000.75
$ws.Range("K103").Value = 3002.25
$ws.Range("M103").Value = -2416.25

$ws.Range("H106").Value = 2749.5
$ws.Range("J106").Value = 2749.5
$ws.Range("L106").Value = 2749.5
$ws.Range("N106").Value = -4011.5

$ws.Range("H127").Value = 402
$ws.Range("I127").Value = 402
$ws.Range("K127").Value = 1206
$ws.Range("M127").Value = 3754

$ws.Range("H129").Value = 53572390
$ws.Range("J129").Value = 1998
$ws.Range("L129").Value = 5994
$ws.Range("N129").Value = -15994

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 2500
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()

$ws.Range("H32").Value = 8954.521000000001
$ws.Range("I32").Value = 4830.4043
$ws.Range("J32").Value = 17030.916
$ws.Range("K32").Value = 4830.4043
$ws.Range("L32").Value = 17030.916
$ws.Range("M32").Value = -4543.4043
$ws.Range("N32").Value = -17604.916

$ws.Range("H44").Value = 8024.8945
$ws.Range("J44").Value = 8248.5
$ws.Range("L44").Value = 8248.5
$ws.Range("N44").Value = -9224.5

$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H61").Value = 2610.5293
$ws.Range("I61").Value = 2292.2666
$ws.Range("K61").Value = 2292.2666
$ws.Range("M61").Value = -2080.2666

$ws.Range("H63").Value = 4950.769
$ws.Range("J63").Value = 5799.5
$ws.Range("L63").Value = 5799.5
$ws.Range("N63").Value = -7171.5

$ws.Range("H66").Value = 4950.769
$ws.Range("J66").Value = 5799.5
$ws.Range("L66").Value = 28997.5
$ws.Range("N66").Value = -35861.5

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H97").Value = 690895.9399999999
$ws.Range("I97").Value = 983468
$ws.Range("J97").Value = 1261.7858
$ws.Range("K97").Value = 983468
$ws.Range("L97").Value = 1261.7858
$ws.Range("M97").Value = -982972
$ws.Range("N97").Value = -2253.7858

$ws.Range("H102").Value = 10420637
$ws.Range("I102").Value = 11907870
$ws.Range("K102").Value = 11907870
$ws.Range("M102").Value = -11906248

$ws.Range("H132").Value = 2529.6428
$ws.Range("I132").Value = 2416.5386
$ws.Range("K132").Value = 7249.6158
$ws.Range("M132").Value = -4719.6158

$ws.Range("H136").Value = 2610.5293
$ws.Range("I136").Value = 2292.2666
$ws.Range("K136").Value = 6876.7998
$ws.Range("M136").Value = -4326.7998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H35").Value = 45999.5
$ws.Range("J35").Value = 45999.5
$ws.Range("L35").Value = 45999.5
$ws.Range("N35").Value = -46619.5

$ws.Range("H82").Value = 30202.5
$ws.Range("I82").Value = 30202.5
$ws.Range("K82").Value = 30202.5
$ws.Range("M82").Value = -29819.5

$ws.Range("H85").Value = 30202.5
$ws.Range("I85").Value = 30202.5
$ws.Range("K85").Value = 30202.5
$ws.Range("M85").Value = -28876.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 116801.75
$ws.Range("I31").Value = 10012
$ws.Range("K31").Value = 10012
$ws.Range("M31").Value = -9717

$ws.Range("H34").Value = 116801.75
$ws.Range("I34").Value = 10012
$ws.Range("K34").Value = 10012
$ws.Range("M34").Value = -9810

$ws.Range("H105").Value = 6244.75
$ws.Range("I105").Value = 6492.5
$ws.Range("J105").Value = 5997
$ws.Range("K105").Value = 6492.5
$ws.Range("L105").Value = 5997
$ws.Range("M105").Value = -4745.5
$ws.Range("N105").Value = -9491

$ws.Range("H122").Value = 2139.9092
$ws.Range("J122").Value = 3252.25
$ws.Range("L122").Value = 9756.75
$ws.Range("N122").Value = -14656.75

$ws.Range("H131").Value = 56448.844
$ws.Range("J131").Value = 56448.844
$ws.Range("L131").Value = 56448.844
$ws.Range("N131").Value = -66528.844

$ws.Range("H134").Value = 37839.04
$ws.Range("I134").Value = 58624.875
$ws.Range("J134").Value = 4581.7
$ws.Range("K134").Value = 175874.625
$ws.Range("L134").Value = 13745.1
$ws.Range("M134").Value = -173339.625
$ws.Range("N134").Value = -18815.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 119.333336
$ws.Range("I6").Value = 119.333336
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 358.000008
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -245.000008
$ws.Range("N6").ClearContents()

$ws.Range("H26").Value = 490.55554
$ws.Range("J26").Value = 999
$ws.Range("L26").Value = 2997
$ws.Range("N26").Value = -3573

$ws.Range("H55").Value = 37134.83
$ws.Range("I55").Value = 358.85715
$ws.Range("J55").Value = 48836.273
$ws.Range("K55").Value = 1076.57145
$ws.Range("L55").Value = 146508.819
$ws.Range("M55").Value = -899.5714499999999
$ws.Range("N55").Value = -146862.819

$ws.Range("H109").Value = 2516.125
$ws.Range("I109").Value = 2705.4
$ws.Range("J109").Value = 2200.6667
$ws.Range("K109").Value = 8116.200000000001
$ws.Range("L109").Value = 6602.000100000001
$ws.Range("M109").Value = -7076.200000000001
$ws.Range("N109").Value = -8682.000100000001

$ws.Range("H113").Value = 4278.684
$ws.Range("I113").Value = 16333.333
$ws.Range("J113").Value = 2018.4375
$ws.Range("K113").Value = 48999.999
$ws.Range("L113").Value = 6055.3125
$ws.Range("M113").Value = -46829.999
$ws.Range("N113").Value = -10395.3125

$ws.Range("H119").Value = 6041.4
$ws.Range("J119").Value = 4766
$ws.Range("L119").Value = 14298
$ws.Range("N119").Value = -23974

$ws.Range("H131").Value = 11264230
$ws.Range("J131").Value = 15876958
$ws.Range("L131").Value = 47630874
$ws.Range("N131").Value = -47640954

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 280073.6
$ws.Range("I122").Value = 357800.56
$ws.Range("K122").Value = 1073401.68
$ws.Range("M122").Value = -1070951.68

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 90037.89999999999
$ws.Range("J22").Value = 1276.7778
$ws.Range("L22").Value = 1276.7778
$ws.Range("N22").Value = -1866.7778

$ws.Range("H27").Value = 90037.89999999999
$ws.Range("J27").Value = 1276.7778
$ws.Range("L27").Value = 1276.7778
$ws.Range("N27").Value = -1490.7778

$ws.Range("H40").Value = 4091.9285
$ws.Range("I40").Value = 2274.0833
$ws.Range("K40").Value = 2274.0833
$ws.Range("M40").Value = -2138.0833

$ws.Range("H51").Value = 41599.8
$ws.Range("J51").Value = 41599.8
$ws.Range("L51").Value = 41599.8
$ws.Range("N51").Value = -42555.8

$ws.Range("H132").Value = 4533.485
$ws.Range("I132").Value = 3903.0476
$ws.Range("K132").Value = 11709.1428
$ws.Range("M132").Value = -9179.1428

$ws.Range("H136").Value = 56844.473
$ws.Range("I136").Value = 89625.74000000001
$ws.Range("K136").Value = 268877.22
$ws.Range("M136").Value = -266327.22

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 21977772
$ws.Range("J132").Value = 2179376.2
$ws.Range("L132").Value = 6538128.600000001
$ws.Range("N132").Value = -6543188.600000001

$ws.Range("H136").Value = 1673.2449
$ws.Range("I136").Value = 1170.2046
$ws.Range("K136").Value = 3510.6138
$ws.Range("M136").Value = -960.6138000000001

